$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

$newRow = 85

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $newRow 1 "7112"
Set-TextCell $newRow 2 "9/1/2025"
Set-TextCell $newRow 3 "OLLEROS 2488"
Set-TextCell $newRow 4 "13"
Set-TextCell $newRow 5 "809371829"
Set-TextCell $newRow 6 "PEBCOM"
Set-TextCell $newRow 7 "Pendiente"
Set-TextCell $newRow 8 "Cambiar "

$ws.Cells.Item($newRow, 9).Value = 1

Set-TextCell $newRow 10 "Cambio"
Set-TextCell $newRow 11 "Sin equipos"
Set-TextCell $newRow 12 "Terminal"

$ws.Cells.Item($newRow, 13).Value = -58.444214
$ws.Cells.Item($newRow, 14).Value = -34.571197

Set-TextCell $newRow 15 "Palermo"
Set-TextCell $newRow 16 "Capital Sur"
